$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    # Prices in column D are plain numeric-looking strings stored as
    # TEXT in the source sheet (t="inlineStr"); a bare .Value assignment
    # would let Excel auto-convert them to Number and silently mangle
    # formatting (e.g. "1.80" -> 1.8). Force text via NumberFormat "@",
    # write the value, then restore the original style so no spurious
    # formatting diff is left behind.
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $newValue
    $range.Style = $origStyle
}

Set-TextValue "D2" "45.332.21"
$ws.Range("E2").Value = "  +2.71%  "
Set-TextValue "D3" "2.423.50"
$ws.Range("E3").Value = "  -0.23%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "317.95"
$ws.Range("E5").Value = "  +3.16%  "
Set-TextValue "D6" "102.59"
$ws.Range("E6").Value = "  +2.20%  "
Set-TextValue "D7" "0.516"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  -0.04%  "
Set-TextValue "D9" "0.528"
$ws.Range("E9").Value = "  +6.00%  "
Set-TextValue "D10" "35.56"
$ws.Range("E10").Value = "  +0.45%  "
Set-TextValue "D11" "0.0802"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  -2.27%  "
Set-TextValue "D13" "18.19"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("E14").Value = "  +2.13%  "
Set-TextValue "D15" "2.804.10"
$ws.Range("E15").Value = "  -0.10%  "
Set-TextValue "D16" "2.437.23"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("E17").Value = "  +1.24%  "
Set-TextValue "D18" "45.262.49"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("E20").Value = "  -1.63%  "
Set-TextValue "D21" "0.0₃0922"
$ws.Range("E21").Value = "  +1.82%  "
Set-TextValue "D22" "68.91"
$ws.Range("E22").Value = "  +0.43%  "
Set-TextValue "D23" "244.36"
$ws.Range("E23").Value = "  +1.72%  "
Set-TextValue "D24" "2.28"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  -0.03%  "
Set-TextValue "D27" "25.64"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("E29").Value = "  -0.05%  "
Set-TextValue "D30" "49.35"
$ws.Range("E30").Value = "  +2.77%  "
Set-TextValue "D31" "32.94"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D32" "0.125"
$ws.Range("E32").Value = "  +5.10%  "
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D33" "19.99"
$ws.Range("E33").Value = "  +6.98%  "
Set-TextValue "D34" "5.22"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  +0.23%  "
Set-TextValue "D36" "0.0762"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D37" "1.87"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D38" "4.46"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("E39").Value = "  -1.79%  "
Set-TextValue "D40" "125.65"
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("E41").Value = "  +0.66%  "
Set-TextValue "D42" "2.21"
$ws.Range("E42").Value = "  -3.52%  "
Set-TextValue "D43" "20.58"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("E44").Value = "  +0.66%  "
Set-TextValue "D45" "1.937.64"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("E47").Value = "  +1.39%  "
Set-TextValue "D48" "1.80"
$ws.Range("E48").Value = "  +8.89%  "
$ws.Range("E49").Value = "  -3.18%  "
Set-TextValue "D50" "76.99"
$ws.Range("E50").Value = "  +4.64%  "
Set-TextValue "D51" "4.77"
$ws.Range("E51").Value = "  +5.25%  "
